$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "45.857.11"
Set-TextValue "E2" "  -0.27%  "
Set-TextValue "D3" "2.577.62"
Set-TextValue "E3" "  +8.42%  "
Set-TextValue "D4" "0.999"
Set-TextValue "E4" "  -0.08%  "
Set-TextValue "D5" "304.23"
Set-TextValue "E5" "  +1.15%  "
Set-TextValue "D6" "98.87"
Set-TextValue "E6" "  +0.44%  "
Set-TextValue "E7" "  +4.56%  "
Set-TextValue "E8" "  +0.00%  "
Set-TextValue "D9" "0.562"
Set-TextValue "E9" "  +10.78%  "
Set-TextValue "D10" "37.72"
Set-TextValue "E10" "  +9.38%  "
Set-TextValue "D11" "0.0827"
Set-TextValue "E11" "  +4.87%  "
Set-TextValue "D12" "7.89"
Set-TextValue "E12" "  +10.84%  "
Set-TextValue "D13" "2.970.39"
Set-TextValue "E13" "  +8.39%  "
Set-TextValue "E14" "  +1.01%  "
Set-TextValue "D15" "2.593.10"
Set-TextValue "E15" "  +9.03%  "
Set-TextValue "D16" "0.886"
Set-TextValue "E16" "  +7.48%  "
Set-TextValue "D17" "14.63"
Set-TextValue "E17" "  +6.55%  "
Set-TextValue "D18" "45.921.41"
Set-TextValue "E18" "  +0.08%  "
Set-TextValue "E19" "  +1.29%  "
Set-TextValue "D20" "0.0₃0994"
Set-TextValue "E20" "  +4.33%  "
Set-TextValue "D21" "6.57"
Set-TextValue "E21" "  +8.61%  "
Set-TextValue "D22" "70.14"
Set-TextValue "E22" "  +5.05%  "
Set-TextValue "D23" "251.33"
Set-TextValue "E23" "  +3.37%  "
Set-TextValue "D24" "2.95"
Set-TextValue "E24" "  +6.09%  "
Set-TextValue "E25" "  +13.13%  "
Set-TextValue "D26" "27.27"
Set-TextValue "E26" "  +30.40%  "
Set-TextValue "D27" "0.999"
Set-TextValue "E27" "  +0.05%  "
Set-TextValue "D28" "10.27"
Set-TextValue "E28" "  +5.45%  "
Set-TextValue "E29" "  +2.27%  "
Set-TextValue "D30" "38.83"
Set-TextValue "E30" "  -2.25%  "
Set-TextValue "D31" "6.02"
Set-TextValue "E31" "  +8.68%  "
Set-TextValue "E32" "  -3.10%  "
Set-TextValue "E33" "  +4.32%  "
Set-TextValue "D34" "2.25"
Set-TextValue "E34" "  +17.33%  "
Set-TextValue "D35" "151.76"
Set-TextValue "E35" "  +3.33%  "
Set-TextValue "D36" "0.0819"
Set-TextValue "E36" "  +5.88%  "
Set-TextValue "E37" "  +1.67%  "
Set-TextValue "E38" "  +4.20%  "
Set-TextValue "D39" "4.12"
Set-TextValue "E39" "  +6.09%  "
Set-TextValue "D40" "15.48"
Set-TextValue "E40" "  +2.55%  "
Set-TextValue "D41" "3.52"
Set-TextValue "E41" "  +9.72%  "
Set-TextValue "D42" "0.0317"
Set-TextValue "E42" "  +6.05%  "
Set-TextValue "D43" "2.040.31"
Set-TextValue "E43" "  +5.52%  "
Set-TextValue "D44" "19.43"
Set-TextValue "E44" "  +37.77%  "
Set-TextValue "D45" "0.999"
Set-TextValue "E45" "  -0.05%  "
Set-TextValue "D46" "90.38"
Set-TextValue "E46" "  -1.75%  "
Set-TextValue "E47" "  +7.91%  "
Set-TextValue "D48" "108.03"
Set-TextValue "E48" "  +9.19%  "
Set-TextValue "E49" "  -2.12%  "
Set-TextValue "D50" "2.830.59"
Set-TextValue "E50" "  +8.46%  "
Set-TextValue "E51" "  +5.62%  "
